$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Description" column text for a few rows: the folder-path
# descriptions used to be phrased "relative to the 'Code' folder" and are
# now phrased "relative to the location of this file".
$ws.Range("C2").Value = 'Path to the folder with pkml simulation files; relative to the location of this file'
$ws.Range("C3").Value = 'Path to the folder with excel files with parametrization; relative to the location of this file'
$ws.Range("C14").Value = 'Path to the folder where the results should be saved to; relative to the location of this file'
$ws.Range("C10").Value = 'Path to the folder where experimental data files are located; relative to the location of this file'

# The default/example compoundPropertiesFile value is no longer prefilled.
$ws.Range("B13").ClearContents()

# Match the author's final selection position in the saved file.
$ws.Range("C10").Select()
